$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- HUD UI: add minimap -> new "HabitatID" column (C) on the Scene sheet ---

# Header row (row 1-3 already carry row-level formatting via s="1" on the
# row, so a plain value write on column C naturally inherits it).
$ws.Range("C1").Value = "HabitatID"
$ws.Range("C2").Value = "int"
$ws.Range("C3").Value = "栖息地 ID"

# Data rows: copy column A's number formatting onto C5:C7 (xlPasteFormats)
# so the new values land as real numbers using the existing style, then
# fill in the values.
$ws.Range("A5").Copy()
$ws.Range("C5:C7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 3

# Restore the active selection to the new column's header cell.
$ws.Range("C4").Select()
